# Update shooter calibration data for Kentwood.
# Rewrites the ballistic table (distance/azimuth/elevation/speed/hopper)
# in rows 2-110, shifting existing values and appending new rows up to A110=180.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 109,5
$data[0,0] = 72
$data[0,1] = 0
$data[0,2] = 1428.02
$data[0,3] = 437.23645671509507
$data[0,4] = 8
$data[1,0] = 73
$data[1,1] = 0
$data[1,2] = 1446.877
$data[1,3] = 437.93711585281255
$data[1,4] = 8
$data[2,0] = 74
$data[2,1] = 0
$data[2,2] = 1465.7339999999999
$data[2,3] = 438.6909663204865
$data[2,4] = 8
$data[3,0] = 75
$data[3,1] = 0
$data[3,2] = 1484.5909999999999
$data[3,3] = 439.49522807209462
$data[3,4] = 8
$data[4,0] = 76
$data[4,1] = 0
$data[4,2] = 1503.4480000000001
$data[4,3] = 440.34717763552874
$data[4,4] = 8
$data[5,0] = 77
$data[5,1] = 0
$data[5,2] = 1522.3050000000001
$data[5,3] = 441.24414811259408
$data[5,4] = 8
$data[6,0] = 78
$data[6,1] = 0
$data[6,2] = 1541.162
$data[6,3] = 442.18352917901086
$data[6,4] = 8
$data[7,0] = 79
$data[7,1] = 0
$data[7,2] = 1560.019
$data[7,3] = 443.16276708441217
$data[7,4] = 8
$data[8,0] = 80
$data[8,1] = 0
$data[8,2] = 1578.876
$data[8,3] = 444.17936465234629
$data[8,4] = 8
$data[9,0] = 81
$data[9,1] = 0
$data[9,2] = 1597.7329999999999
$data[9,3] = 445.23088128027496
$data[9,4] = 8
$data[10,0] = 82
$data[10,1] = 0
$data[10,2] = 1616.59
$data[10,3] = 446.31493293957374
$data[10,4] = 8
$data[11,0] = 83
$data[11,1] = 0
$data[11,2] = 1635.4469999999999
$data[11,3] = 447.4291921755331
$data[11,4] = 8
$data[12,0] = 84
$data[12,1] = 0
$data[12,2] = 1654.3039999999999
$data[12,3] = 448.57138810735626
$data[12,4] = 8
$data[13,0] = 85
$data[13,1] = 0
$data[13,2] = 1673.1610000000001
$data[13,3] = 449.73930642816163
$data[13,4] = 8
$data[14,0] = 86
$data[14,1] = 0
$data[14,2] = 1692.018
$data[14,3] = 450.93078940498094
$data[14,4] = 8
$data[15,0] = 87
$data[15,1] = 0
$data[15,2] = 1710.875
$data[15,3] = 452.14373587876048
$data[15,4] = 8
$data[16,0] = 88
$data[16,1] = 0
$data[16,2] = 1729.732
$data[16,3] = 453.37610126436027
$data[16,4] = 8
$data[17,0] = 89
$data[17,1] = 0
$data[17,2] = 1748.5889999999999
$data[17,3] = 454.62589755055444
$data[17,4] = 8
$data[18,0] = 90
$data[18,1] = 0
$data[18,2] = 1767.4459999999999
$data[18,3] = 455.89119330003086
$data[18,4] = 8
$data[19,0] = 91
$data[19,1] = 0
$data[19,2] = 1786.3029999999999
$data[19,3] = 457.17011364939196
$data[19,4] = 8
$data[20,0] = 92
$data[20,1] = 0
$data[20,2] = 1805.16
$data[20,3] = 458.46084030915415
$data[20,4] = 8
$data[21,0] = 93
$data[21,1] = 0
$data[21,2] = 1824.0170000000001
$data[21,3] = 459.76161156374712
$data[21,4] = 8
$data[22,0] = 94
$data[22,1] = 0
$data[22,2] = 1842.874
$data[22,3] = 461.07072227151582
$data[22,4] = 8
$data[23,0] = 95
$data[23,1] = 0
$data[23,2] = 1861.731
$data[23,3] = 462.38652386471847
$data[23,4] = 8
$data[24,0] = 96
$data[24,1] = 0
$data[24,2] = 1880.588
$data[24,3] = 463.70742434952706
$data[24,4] = 8
$data[25,0] = 97
$data[25,1] = 0
$data[25,2] = 1899.4449999999999
$data[25,3] = 465.03188830602869
$data[25,4] = 8
$data[26,0] = 98
$data[26,1] = 0
$data[26,2] = 1918.3019999999999
$data[26,3] = 466.35843688822308
$data[26,4] = 8
$data[27,0] = 99
$data[27,1] = 0
$data[27,2] = 1937.1589999999999
$data[27,3] = 467.6856478240253
$data[27,4] = 8
$data[28,0] = 100
$data[28,1] = 0
$data[28,2] = 1956.0159999999998
$data[28,3] = 469.01215541526369
$data[28,4] = 8
$data[29,0] = 101
$data[29,1] = 0
$data[29,2] = 1974.873
$data[29,3] = 470.33665053768084
$data[29,4] = 8
$data[30,0] = 102
$data[30,1] = 0
$data[30,2] = 1993.73
$data[30,3] = 471.65788064093329
$data[30,4] = 8
$data[31,0] = 103
$data[31,1] = 0
$data[31,2] = 2012.587
$data[31,3] = 472.97464974859179
$data[31,4] = 8
$data[32,0] = 104
$data[32,1] = 0
$data[32,2] = 2031.444
$data[32,3] = 474.28581845814108
$data[32,4] = 8
$data[33,0] = 105
$data[33,1] = 0
$data[33,2] = 2050.3009999999999
$data[33,3] = 475.59030394098011
$data[33,4] = 8
$data[34,0] = 106
$data[34,1] = 0
$data[34,2] = 2069.1579999999999
$data[34,3] = 476.88707994242088
$data[34,4] = 8
$data[35,0] = 107
$data[35,1] = 0
$data[35,2] = 2088.0149999999999
$data[35,3] = 478.17517678169077
$data[35,4] = 8
$data[36,0] = 108
$data[36,1] = 0
$data[36,2] = 2106.8719999999998
$data[36,3] = 479.45368135193064
$data[36,4] = 8
$data[37,0] = 109
$data[37,1] = 0
$data[37,2] = 2125.7289999999998
$data[37,3] = 480.7217371201956
$data[37,4] = 8
$data[38,0] = 110
$data[38,1] = 0
$data[38,2] = 2144.5859999999998
$data[38,3] = 481.97854412745403
$data[38,4] = 8
$data[39,0] = 111
$data[39,1] = 0
$data[39,2] = 2163.4429999999998
$data[39,3] = 483.22335898858944
$data[39,4] = 8
$data[40,0] = 112
$data[40,1] = 0
$data[40,2] = 2182.2999999999997
$data[40,3] = 484.45549489239863
$data[40,4] = 8
$data[41,0] = 113
$data[41,1] = 0
$data[41,2] = 2201.1569999999997
$data[41,3] = 485.67432160159194
$data[41,4] = 8
$data[42,0] = 114
$data[42,1] = 0
$data[42,2] = 2220.0139999999997
$data[42,3] = 486.87926545279595
$data[42,4] = 8
$data[43,0] = 115
$data[43,1] = 0
$data[43,2] = 2238.8709999999996
$data[43,3] = 488.06980935654826
$data[43,4] = 8
$data[44,0] = 116
$data[44,1] = 0
$data[44,2] = 2257.7279999999996
$data[44,3] = 489.24549279730365
$data[44,4] = 8
$data[45,0] = 117
$data[45,1] = 0
$data[45,2] = 2276.5849999999996
$data[45,3] = 490.40591183342792
$data[45,4] = 8
$data[46,0] = 118
$data[46,1] = 0
$data[46,2] = 2295.4419999999996
$data[46,3] = 491.55071909720243
$data[46,4] = 8
$data[47,0] = 119
$data[47,1] = 0
$data[47,2] = 2314.2989999999995
$data[47,3] = 492.67962379482344
$data[47,4] = 8
$data[48,0] = 120
$data[48,1] = 0
$data[48,2] = 2333.1559999999999
$data[48,3] = 493.79239170639937
$data[48,4] = 8
$data[49,0] = 121
$data[49,1] = 0
$data[49,2] = 2352.0129999999999
$data[49,3] = 494.88884518595376
$data[49,4] = 8
$data[50,0] = 122
$data[50,1] = 0
$data[50,2] = 2370.87
$data[50,3] = 495.96886316142434
$data[50,4] = 8
$data[51,0] = 123
$data[51,1] = 0
$data[51,2] = 2389.7269999999999
$data[51,3] = 497.03238113466193
$data[51,4] = 8
$data[52,0] = 124
$data[52,1] = 0
$data[52,2] = 2408.5839999999998
$data[52,3] = 498.07939118143315
$data[52,4] = 8
$data[53,0] = 125
$data[53,1] = 0
$data[53,2] = 2427.4409999999998
$data[53,3] = 499.10994195141632
$data[53,4] = 8
$data[54,0] = 126
$data[54,1] = 0
$data[54,2] = 2446.2979999999998
$data[54,3] = 500.12413866820577
$data[54,4] = 8
$data[55,0] = 127
$data[55,1] = 0
$data[55,2] = 2465.1549999999997
$data[55,3] = 501.12214312930848
$data[55,4] = 8
$data[56,0] = 128
$data[56,1] = 0
$data[56,2] = 2484.0119999999997
$data[56,3] = 502.1041737061463
$data[56,4] = 8
$data[57,0] = 129
$data[57,1] = 0
$data[57,2] = 2502.8689999999997
$data[57,3] = 503.07050534405528
$data[57,4] = 8
$data[58,0] = 130
$data[58,1] = 0
$data[58,2] = 2521.7259999999997
$data[58,3] = 504.02146956228432
$data[58,4] = 8
$data[59,0] = 131
$data[59,1] = 0
$data[59,2] = 2540.5829999999996
$data[59,3] = 504.95745445399837
$data[59,4] = 8
$data[60,0] = 132
$data[60,1] = 0
$data[60,2] = 2559.4399999999996
$data[60,3] = 505.87890468627472
$data[60,4] = 8
$data[61,0] = 133
$data[61,1] = 0
$data[61,2] = 2578.2969999999996
$data[61,3] = 506.78632150010424
$data[61,4] = 8
$data[62,0] = 134
$data[62,1] = 0
$data[62,2] = 2597.1539999999995
$data[62,3] = 507.68026271039378
$data[62,4] = 8
$data[63,0] = 135
$data[63,1] = 0
$data[63,2] = 2616.0109999999995
$data[63,3] = 508.56134270596306
$data[63,4] = 8
$data[64,0] = 136
$data[64,1] = 0
$data[64,2] = 2634.8679999999995
$data[64,3] = 509.43023244954577
$data[64,4] = 8
$data[65,0] = 137
$data[65,1] = 0
$data[65,2] = 2653.7249999999999
$data[65,3] = 510.28765947779073
$data[65,4] = 8
$data[66,0] = 138
$data[66,1] = 0
$data[66,2] = 2672.5819999999999
$data[66,3] = 511.1344079012589
$data[66,4] = 8
$data[67,0] = 139
$data[67,1] = 0
$data[67,2] = 2691.4389999999999
$data[67,3] = 511.97131840442705
$data[67,4] = 8
$data[68,0] = 140
$data[68,1] = 0
$data[68,2] = 2710.2959999999998
$data[68,3] = 512.79928824568435
$data[68,4] = 8
$data[69,0] = 141
$data[69,1] = 0
$data[69,2] = 2729.1529999999998
$data[69,3] = 513.61927125733644
$data[69,4] = 8
$data[70,0] = 142
$data[70,1] = 0
$data[70,2] = 2748.0099999999998
$data[70,3] = 514.43227784559997
$data[70,4] = 8
$data[71,0] = 143
$data[71,1] = 0
$data[71,2] = 2766.8669999999997
$data[71,3] = 515.23937499060901
$data[71,4] = 8
$data[72,0] = 144
$data[72,1] = 0
$data[72,2] = 2785.7239999999997
$data[72,3] = 516.04168624640681
$data[72,4] = 8
$data[73,0] = 145
$data[73,1] = 0
$data[73,2] = 2804.5809999999997
$data[73,3] = 516.84039174095608
$data[73,4] = 8
$data[74,0] = 146
$data[74,1] = 0
$data[74,2] = 2823.4379999999996
$data[74,3] = 517.6367281761319
$data[74,4] = 8
$data[75,0] = 147
$data[75,1] = 0
$data[75,2] = 2842.2949999999996
$data[75,3] = 518.43198882772015
$data[75,4] = 8
$data[76,0] = 148
$data[76,1] = 0
$data[76,2] = 2861.1519999999996
$data[76,3] = 519.22752354542502
$data[76,4] = 8
$data[77,0] = 149
$data[77,1] = 0
$data[77,2] = 2880.0089999999996
$data[77,3] = 520.02473875286171
$data[77,4] = 8
$data[78,0] = 150
$data[78,1] = 0
$data[78,2] = 2898.8659999999995
$data[78,3] = 520.82509744756146
$data[78,4] = 8
$data[79,0] = 151
$data[79,1] = 0
$data[79,2] = 2917.7229999999995
$data[79,3] = 521.63011920096835
$data[79,4] = 8
$data[80,0] = 152
$data[80,1] = 0
$data[80,2] = 2936.58
$data[80,3] = 522.44138015844248
$data[80,4] = 8
$data[81,0] = 153
$data[81,1] = 0
$data[81,2] = 2955.4369999999999
$data[81,3] = 523.26051303925362
$data[81,4] = 8
$data[82,0] = 154
$data[82,1] = 0
$data[82,2] = 2974.2939999999999
$data[82,3] = 524.08920713659074
$data[82,4] = 8
$data[83,0] = 155
$data[83,1] = 0
$data[83,2] = 2993.1509999999998
$data[83,3] = 524.92920831755475
$data[83,4] = 8
$data[84,0] = 156
$data[84,1] = 0
$data[84,2] = 3012.0079999999998
$data[84,3] = 525.78231902315804
$data[84,4] = 8
$data[85,0] = 157
$data[85,1] = 0
$data[85,2] = 3030.8649999999998
$data[85,3] = 526.65039826833174
$data[85,4] = 8
$data[86,0] = 158
$data[86,1] = 0
$data[86,2] = 3049.7219999999998
$data[86,3] = 527.53536164191667
$data[86,4] = 8
$data[87,0] = 159
$data[87,1] = 0
$data[87,2] = 3068.5789999999997
$data[87,3] = 528.439181306671
$data[87,4] = 8
$data[88,0] = 160
$data[88,1] = 0
$data[88,2] = 3087.4359999999997
$data[88,3] = 529.36388599926488
$data[88,4] = 8
$data[89,0] = 161
$data[89,1] = 0
$data[89,2] = 3106.2929999999997
$data[89,3] = 530.31156103028582
$data[89,4] = 8
$data[90,0] = 162
$data[90,1] = 0
$data[90,2] = 3125.1499999999996
$data[90,3] = 531.2843482842278
$data[90,4] = 8
$data[91,0] = 163
$data[91,1] = 0
$data[91,2] = 3144.0069999999996
$data[91,3] = 532.28444621950769
$data[91,4] = 8
$data[92,0] = 164
$data[92,1] = 0
$data[92,2] = 3162.8639999999996
$data[92,3] = 533.31410986845106
$data[92,4] = 8
$data[93,0] = 165
$data[93,1] = 0
$data[93,2] = 3181.7209999999995
$data[93,3] = 534.3756508373009
$data[93,4] = 8
$data[94,0] = 166
$data[94,1] = 0
$data[94,2] = 3200.5779999999995
$data[94,3] = 535.47143730620849
$data[94,4] = 8
$data[95,0] = 167
$data[95,1] = 0
$data[95,2] = 3219.4349999999995
$data[95,3] = 536.60389402924659
$data[95,4] = 8
$data[96,0] = 168
$data[96,1] = 0
$data[96,2] = 3238.2919999999995
$data[96,3] = 537.77550233439626
$data[96,4] = 8
$data[97,0] = 169
$data[97,1] = 0
$data[97,2] = 3257.1489999999999
$data[97,3] = 538.98880012355596
$data[97,4] = 8
$data[98,0] = 170
$data[98,1] = 0
$data[98,2] = 3276.0059999999999
$data[98,3] = 540.24638187253697
$data[98,4] = 8
$data[99,0] = 171
$data[99,1] = 0
$data[99,2] = 3294.8629999999998
$data[99,3] = 541.55089863106252
$data[99,4] = 8
$data[100,0] = 172
$data[100,1] = 0
$data[100,2] = 3313.72
$data[100,3] = 542.9050580227746
$data[100,4] = 8
$data[101,0] = 173
$data[101,1] = 0
$data[101,2] = 3332.5769999999998
$data[101,3] = 544.31162424522574
$data[101,4] = 8
$data[102,0] = 174
$data[102,1] = 0
$data[102,2] = 3351.4339999999997
$data[102,3] = 545.77341806988272
$data[102,4] = 8
$data[103,0] = 175
$data[103,1] = 0
$data[103,2] = 3370.2909999999997
$data[103,3] = 547.29331684212741
$data[103,4] = 8
$data[104,0] = 176
$data[104,1] = 0
$data[104,2] = 3389.1479999999997
$data[104,3] = 548.87425448125589
$data[104,4] = 8
$data[105,0] = 177
$data[105,1] = 0
$data[105,2] = 3408.0049999999997
$data[105,3] = 550.51922148047709
$data[105,4] = 8
$data[106,0] = 178
$data[106,1] = 0
$data[106,2] = 3426.8619999999996
$data[106,3] = 552.23126490691504
$data[106,4] = 8
$data[107,0] = 179
$data[107,1] = 0
$data[107,2] = 3445.7189999999996
$data[107,3] = 554.01348840160529
$data[107,4] = 8
$data[108,0] = 180
$data[108,1] = 0
$data[108,2] = 3464.5759999999996
$data[108,3] = 555.86905217950255
$data[108,4] = 8

$ws.Range("A2:E110").Value = $data

# Select the full used range, matching the post-edit selection state.
$ws.Range("A1:E110").Select()
